$d = $word.ActiveDocument

# The first paragraph of the document (top-right header table cell) currently
# reads "In the Online Civil Claims at Northampton " in bold, sz16. Replace it
# with "In the County Court" in non-bold, sz24 (12pt), matching the template
# update. Scope the Find/Replace to the paragraph's own Range so it collapses
# onto just the replaced text (rather than the whole document) afterwards.
$p = $d.Paragraphs.Item(1)
$r = $p.Range

$found = $r.Find.Execute("In the Online Civil Claims at Northampton ", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "In the County Court", 2)

$r.Font.Bold = $false
$r.Font.Size = 12
